# Auto-generated Excel COM-interop edit script
# Applies the numeric cell updates described by the upstream diff
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 12500
$ws.Range("I47").Value = 12500
$ws.Range("K47").Value = 12500
$ws.Range("M47").Value = -11528

$ws.Range("H54").Value = 33646.668
$ws.Range("I54").Value = 12333.333
$ws.Range("K54").Value = 12333.333
$ws.Range("M54").Value = -11847.333

$ws.Range("H98").Value = 42366.367
$ws.Range("I98").Value = 1417.75
$ws.Range("J98").Value = 260759
$ws.Range("K98").Value = 1417.75
$ws.Range("L98").Value = 260759
$ws.Range("M98").Value = 80.25
$ws.Range("N98").Value = -263755

$ws.Range("H111").Value = 1406.7778
$ws.Range("I111").Value = 1471.5
$ws.Range("J111").Value = 1277.3334
$ws.Range("K111").Value = 4414.5
$ws.Range("L111").Value = 3832.0002
$ws.Range("M111").Value = -1347.5
$ws.Range("N111").Value = -9966.0002

$ws.Range("H113").Value = 2599.8333
$ws.Range("I113").Value = 2519.8
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2519.8
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 734.1999999999998
$ws.Range("N113").Value = -9508

$ws.Range("H116").Value = 7629.56
$ws.Range("I116").Value = 6129.143
$ws.Range("J116").Value = 8213.056
$ws.Range("K116").Value = 6129.143
$ws.Range("L116").Value = 8213.056
$ws.Range("M116").Value = -2687.143
$ws.Range("N116").Value = -15097.056

$ws.Range("H122").Value = 42366.367
$ws.Range("I122").Value = 1417.75
$ws.Range("J122").Value = 260759
$ws.Range("K122").Value = 4253.25
$ws.Range("L122").Value = 782277
$ws.Range("M122").Value = -1803.25
$ws.Range("N122").Value = -787177

$ws.Range("H137").Value = 3216.236
$ws.Range("I137").Value = 1221.0667
$ws.Range("J137").Value = 3741.2808
$ws.Range("K137").Value = 3663.2001
$ws.Range("L137").Value = 11223.8424
$ws.Range("M137").Value = -1113.2001
$ws.Range("N137").Value = -16323.8424

$ws.Range("H138").Value = 2066.4255
$ws.Range("I138").Value = 1576.5897
$ws.Range("J138").Value = 2413.7637
$ws.Range("K138").Value = 4729.7691
$ws.Range("L138").Value = 7241.2911
$ws.Range("M138").Value = 410.2309000000005
$ws.Range("N138").Value = -17521.2911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1758.3
$ws.Range("I122").Value = 1722.7059
$ws.Range("J122").Value = 1960
$ws.Range("K122").Value = 5168.1177
$ws.Range("L122").Value = 5880
$ws.Range("M122").Value = -2718.1177
$ws.Range("N122").Value = -10780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 40780
$ws.Range("J122").Value = 40780
$ws.Range("L122").Value = 40780
$ws.Range("N122").Value = -50580

$ws.Range("H123").Value = 24777.777
$ws.Range("J123").Value = 24777.777
$ws.Range("L123").Value = 24777.777
$ws.Range("N123").Value = -34577.777

$ws.Range("H130").Value = 51885
$ws.Range("J130").Value = 51885
$ws.Range("L130").Value = 51885
$ws.Range("N130").Value = -61925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5425.2104
$ws.Range("I31").Value = 2348.2144
$ws.Range("J31").Value = 6427.0234
$ws.Range("K31").Value = 2348.2144
$ws.Range("L31").Value = 6427.0234
$ws.Range("M31").Value = -2053.2144
$ws.Range("N31").Value = -7017.0234

$ws.Range("H34").Value = 5425.2104
$ws.Range("I34").Value = 2348.2144
$ws.Range("J34").Value = 6427.0234
$ws.Range("K34").Value = 2348.2144
$ws.Range("L34").Value = 6427.0234
$ws.Range("M34").Value = -2146.2144
$ws.Range("N34").Value = -6831.0234

$ws.Range("H80").Value = 33037.75
$ws.Range("J80").Value = 33037.75
$ws.Range("L80").Value = 33037.75
$ws.Range("N80").Value = -35283.75

$ws.Range("H83").Value = 33037.75
$ws.Range("J83").Value = 33037.75
$ws.Range("L83").Value = 99113.25
$ws.Range("N83").Value = -110345.25

$ws.Range("H100").Value = 46996
$ws.Range("J100").Value = 46996
$ws.Range("L100").Value = 46996
$ws.Range("N100").Value = -49160

$ws.Range("H107").Value = 610.8
$ws.Range("I107").Value = 552.2857
$ws.Range("K107").Value = 552.2857
$ws.Range("M107").Value = 1367.7143

$ws.Range("H110").Value = 40798.332
$ws.Range("J110").Value = 40798.332
$ws.Range("L110").Value = 40798.332
$ws.Range("N110").Value = -48978.332

$ws.Range("H132").Value = 56592.23
$ws.Range("I132").Value = 1950.25
$ws.Range("J132").Value = 144019.4
$ws.Range("K132").Value = 5850.75
$ws.Range("L132").Value = 432058.2
$ws.Range("M132").Value = -3320.75
$ws.Range("N132").Value = -437118.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1502.8572
$ws.Range("J7").Value = 116.666664
$ws.Range("L7").Value = 349.999992
$ws.Range("N7").Value = -573.999992

$ws.Range("H34").Value = 1608.8889
$ws.Range("J34").Value = 2011.4286
$ws.Range("L34").Value = 6034.2858
$ws.Range("N34").Value = -6202.2858

$ws.Range("H39").Value = 1175.4375
$ws.Range("I39").Value = 676.75
$ws.Range("J39").Value = 1341.6666
$ws.Range("K39").Value = 2030.25
$ws.Range("L39").Value = 4024.9998
$ws.Range("M39").Value = -1736.25
$ws.Range("N39").Value = -4612.9998

$ws.Range("H55").Value = 1038.5
$ws.Range("I55").Value = 602
$ws.Range("J55").Value = 1475
$ws.Range("K55").Value = 1806
$ws.Range("L55").Value = 4425
$ws.Range("M55").Value = -1629
$ws.Range("N55").Value = -4779

$ws.Range("H129").Value = 1912.25
$ws.Range("I129").Value = 1627.5
$ws.Range("J129").Value = 2007.1666
$ws.Range("K129").Value = 4882.5
$ws.Range("L129").Value = 6021.4998
$ws.Range("M129").Value = 117.5
$ws.Range("N129").Value = -16021.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 12451
$ws.Range("J15").Value = 12451
$ws.Range("L15").Value = 12451
$ws.Range("N15").Value = -13027

$ws.Range("H81").Value = 12451
$ws.Range("J81").Value = 12451
$ws.Range("L81").Value = 12451
$ws.Range("N81").Value = -14447

$ws.Range("H84").Value = 12451
$ws.Range("J84").Value = 12451
$ws.Range("L84").Value = 37353
$ws.Range("N84").Value = -47337

$ws.Range("H122").Value = 2007.3572
$ws.Range("J122").Value = 2101.1428
$ws.Range("L122").Value = 6303.428400000001
$ws.Range("N122").Value = -11203.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2569.842
$ws.Range("I7").Value = 2224.5293
$ws.Range("J7").Value = 5505
$ws.Range("K7").Value = 2224.5293
$ws.Range("L7").Value = 5505
$ws.Range("M7").Value = -2112.5293
$ws.Range("N7").Value = -5729

$ws.Range("H16").Value = 2024.875
$ws.Range("I16").Value = 2024.875
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2024.875
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1854.875
$ws.Range("N16").ClearContents()

$ws.Range("H93").Value = 2224.1667
$ws.Range("I93").Value = 2109
$ws.Range("J93").Value = 2800
$ws.Range("K93").Value = 2109
$ws.Range("L93").Value = 2800
$ws.Range("M93").Value = -861
$ws.Range("N93").Value = -5296

$ws.Range("H112").Value = 43590
$ws.Range("J112").Value = 43590
$ws.Range("L112").Value = 43590
$ws.Range("N112").Value = -46544

$ws.Range("H126").Value = 2569.842
$ws.Range("I126").Value = 2224.5293
$ws.Range("J126").Value = 5505
$ws.Range("K126").Value = 6673.5879
$ws.Range("L126").Value = 16515
$ws.Range("M126").Value = -4203.5879
$ws.Range("N126").Value = -21455

$ws.Range("H133").Value = 27531.273
$ws.Range("J133").Value = 27531.273
$ws.Range("L133").Value = 27531.273
$ws.Range("N133").Value = -32591.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 33750
$ws.Range("J86").Value = 33750
$ws.Range("L86").Value = 33750
$ws.Range("N86").Value = -35996

$ws.Range("H89").Value = 33750
$ws.Range("J89").Value = 33750
$ws.Range("L89").Value = 168750
$ws.Range("N89").Value = -179982

$ws.Range("H109").Value = 32136.445
$ws.Range("J109").Value = 32136.445
$ws.Range("L109").Value = 32136.445
$ws.Range("N109").Value = -34910.445
